$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.172.66"
$ws.Range("E2").Value = "  +1.02%  "
$ws.Range("D3").Value = "3.836.36"
$ws.Range("E3").Value = "  +1.29%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").NumberFormat = "General"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "717.15"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.77%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "172.74"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.75%  "
$ws.Range("D7").Value = "3.835.89"
$ws.Range("E7").Value = "  +1.33%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("E9").Value = "  -0.01%  "
$ws.Range("E10").Value = "  +0.77%  "
$ws.Range("E11").Value = "  +0.80%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.459"
$ws.Range("D12").NumberFormat = "General"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.24%  "
$ws.Range("E13").Value = "  +0.79%  "
$ws.Range("E14").Value = "  +2.05%  "
$ws.Range("D15").Value = "4.482.12"
$ws.Range("E15").Value = "  +1.19%  "
$ws.Range("D16").Value = "3.849.04"
$ws.Range("E16").Value = "  +1.45%  "
$ws.Range("D17").Value = "71.137.74"
$ws.Range("E17").Value = "  +0.92%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.22"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.84%  "
$ws.Range("E19").Value = "  +0.73%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.40"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.42%  "
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.74"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.80%  "
$ws.Range("B22").Value = "BitcoinCash"
$ws.Range("C22").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "494.89"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.50%  "
$ws.Range("E23").Value = "  +2.25%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "85.13"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.76%  "
$ws.Range("E25").Value = "  +3.43%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.68"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.55%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.16"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.88%  "
$ws.Range("E28").Value = "  +3.78%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.11"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.52%  "
$ws.Range("E30").Value = "  -0.03%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.52"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.30%  "
$ws.Range("E32").Value = "  -1.05%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "29.43"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.62%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.181"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.48%  "
$ws.Range("E35").Value = "  +0.17%  "
$ws.Range("D36").Value = "3.800.50"
$ws.Range("E36").Value = "  +1.63%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.07%  "
$ws.Range("E38").Value = "  +0.51%  "
$ws.Range("B39").Value = "Mantle"
$ws.Range("C39").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.03"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.90%  "
$ws.Range("B40").Value = "Filecoin"
$ws.Range("C40").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.01"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.74%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.37"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.08%  "
$ws.Range("E42").Value = "  +2.46%  "
$ws.Range("E43").Value = "  +0.02%  "
$ws.Range("E44").Value = "  +0.04%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.000320"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.04%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "163.46"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.10%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "48.80"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.07%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "422.25"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.76%  "
$ws.Range("E49").Value = "  +1.24%  "
$ws.Range("E50").Value = "  +0.88%  "
$ws.Range("E51").Value = "  -0.57%  "
